$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 65.187973903399524
$ws.Range("C2").Value = 35.761089351678933
$ws.Range("D2").Value = 64.474610650314517
$ws.Range("E2").Value = 39.866832231383704

$ws.Range("B3").Value = 62.557460325816436
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = 51.95158270878715
$ws.Range("E3").Value = 46.752859715113551

$ws.Range("B1:E3").Select()
